$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data_Parking")

# Row 2: update existing record (RCCF584T2B12 -> NJ9E1D0UMB22)
$ws.Range("A2").Value = "NJ9E1D0UMB22"
$ws.Range("B2").Value = "Y 0313 HAS"
$ws.Range("C2").Value = "Motor"
$ws.Range("D2").Value = "2025-01-31 21:09:22"
$ws.Range("E2").Value = "2025-01-31 21:09:47"
$ws.Range("F2").Value = "00:00:25"
$ws.Range("G2").Value = 2000
$ws.Range("H2").Value = 5000
$ws.Range("I2").Value = "Reza Ramdan Permana"
$ws.Range("J2").Value = "./capture/masuk/NJ9E1D0UMB22.png"
$ws.Range("K2").Value = "./capture/keluar/NJ9E1D0UMB22.png"

# Row 3: update existing record (IULMD53YWRNI -> SN12XKBEG18L)
$ws.Range("A3").Value = "SN12XKBEG18L"
$ws.Range("B3").Value = "D 4230 ASQ"
$ws.Range("C3").Value = "Mobil"
$ws.Range("D3").Value = "2025-01-31 21:10:39"
$ws.Range("E3").Value = "2025-01-31 21:11:58"
$ws.Range("F3").Value = "00:01:19"
$ws.Range("G3").Value = 4000
$ws.Range("H3").Value = 60000
$ws.Range("I3").Value = "Reza Ramdan Permana"
$ws.Range("J3").Value = "./capture/masuk/SN12XKBEG18L.png"
$ws.Range("K3").Value = "./capture/keluar/SN12XKBEG18L.png"

# Row 4: new record (58O0J9BUUNER)
$ws.Range("A4").Value = "58O0J9BUUNER"
$ws.Range("B4").Value = "D 9530 JFD"
$ws.Range("C4").Value = "Motor"
$ws.Range("D4").Value = "2025-01-31 21:15:34"
$ws.Range("E4").Value = "2025-01-31 21:16:09"
$ws.Range("F4").Value = "00:00:35"
$ws.Range("G4").Value = 2000
$ws.Range("H4").Value = 5000
$ws.Range("I4").Value = "Reza Ramdan Permana"
$ws.Range("J4").Value = "./capture/masuk/58O0J9BUUNER.png"
$ws.Range("K4").Value = "./capture/keluar/58O0J9BUUNER.png"
